$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order per row: B (Coin), C (Link), D (Price, forced text), E (Volume)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.040.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.756.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("E8").Value = "  -2.73%  "

# Row 9
$ws.Range("E9").Value = "  -1.64%  "

# Row 10
$ws.Range("E10").Value = "  +3.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -14.36%  "

# Row 12
$ws.Range("E12").Value = "  -1.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.243.41"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.99"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.65%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.666.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.29%  "

# Row 16
$ws.Range("E16").Value = "  -2.68%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.760.69"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("E19").Value = "  -1.97%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.96%  "

# Row 22
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.533"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.56%  "

# Row 25
$ws.Range("E25").Value = "  -0.92%  "

# Row 26
$ws.Range("E26").Value = "  +0.36%  "

# Row 27
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0909"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.00%  "

# Row 31
$ws.Range("E31").Value = "  -0.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.93"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.96%  "

# Row 33
$ws.Range("E33").Value = "  -0.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.01%  "

# Row 38
$ws.Range("E38").Value = "  -1.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "350.12"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.77%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.19"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.45%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0589"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.73"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.41%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0255"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.633"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.94%  "

# Row 49
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.19%  "

